# Week 13 logging update for Saints "Players Data.xlsx"
# - Rushing sheet: row 2/3 player names swapped (T.Siemian <-> T.Hill) and
#   their week-13 counting stats corrected; a few other rows' stats bumped.
# - Receiving sheet: several rows' week-13 counting stats corrected.

$wb = $excel.ActiveWorkbook

# ---- Rushing ----
$ws = $wb.Worksheets.Item("Rushing")

# Row 2 / Row 3 players were swapped (the name shown for each row flips)
$ws.Range("B2").Value = "T.Hill"
$ws.Range("B3").Value = "T.Siemian"

# Row 2 stats (now T.Hill)
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 11

# Row 3 stats (now T.Siemian)
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 3

# Other rows with updated stats
$ws.Range("C9").Value = 21
$ws.Range("D9").Value = 17
$ws.Range("C10").Value = 4
$ws.Range("D11").Value = 7

# ---- Receiving ----
$ws = $wb.Worksheets.Item("Receiving")

$ws.Range("C7").Value = 15
$ws.Range("D7").Value = 11

$ws.Range("C9").Value = 27
$ws.Range("D9").Value = 18
$ws.Range("E9").Value = 11

$ws.Range("C10").Value = 35
$ws.Range("D10").Value = 23
$ws.Range("E10").Value = 19

$ws.Range("C11").Value = 40
$ws.Range("D11").Value = 26
$ws.Range("E11").Value = 12
$ws.Range("F11").Value = 5

$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 4

$ws.Range("C14").Value = 21
$ws.Range("D14").Value = 14

$ws.Range("C16").Value = 5
$ws.Range("E16").Value = 5

$ws.Range("E18").Value = 7
$ws.Range("F18").Value = 3

$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 4

$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 2
